$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.6122626666666666
$ws.Range("H2").Value = 1.836788
$ws.Range("I2").Value = 0.006779070576782467
$ws.Range("J2").Value = 0.006779070576782467
$ws.Range("M2").Value = 1.819857
$ws.Range("N2").Value = 5.459571
$ws.Range("O2").Value = 0.01485317462584607
$ws.Range("P2").Value = 0.01485317462584607
$ws.Range("Q2").Value = 1.114230499772
$ws.Range("R2").Value = 10.028074497948
$ws.Range("S2").Value = 0.000100690719077885
$ws.Range("T2").Value = 0.000100690719077885
# Row 3
$ws.Range("G3").Value = 0.6122626666666666
$ws.Range("H3").Value = 1.836788
$ws.Range("I3").Value = 0.006779070576782467
$ws.Range("J3").Value = 0.006779070576782467
$ws.Range("O3").Value = 0.726618572334523
$ws.Range("P3").Value = 0.7266185723345231
$ws.Range("Q3").Value = 54.50825129242666
$ws.Range("R3").Value = 490.5742616318399
$ws.Range("S3").Value = 0.004925798584256648
$ws.Range("T3").Value = 0.004925798584256649
# Row 4
$ws.Range("G4").Value = 0.6122626666666666
$ws.Range("H4").Value = 1.836788
$ws.Range("I4").Value = 0.006779070576782467
$ws.Range("J4").Value = 0.006779070576782467
$ws.Range("M4").Value = 31.52924033333333
$ws.Range("N4").Value = 94.58772099999999
$ws.Range("O4").Value = 0.257333028084772
$ws.Range("P4").Value = 0.257333028084772
$ws.Range("Q4").Value = 19.30417676446088
$ws.Range("R4").Value = 173.737590880148
$ws.Range("S4").Value = 0.001744478759123814
$ws.Range("T4").Value = 0.001744478759123814
# Row 5
$ws.Range("G5").Value = 0.6122626666666666
$ws.Range("H5").Value = 1.836788
$ws.Range("I5").Value = 0.006779070576782467
$ws.Range("J5").Value = 0.006779070576782467
$ws.Range("K5").Value = 1.0
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1464426666666667
$ws.Range("N5").Value = 0.439328
$ws.Range("O5").Value = 0.001195224954858853
$ws.Range("P5").Value = 0.001195224954858853
$ws.Range("Q5").Value = 0.0896613776071111
$ws.Range("R5").Value = 0.8069523984639999
$ws.Range("S5").Value = 0.000008102514324119799
$ws.Range("T5").Value = 0.0000081025143241198
# Row 6
$ws.Range("I6").Value = 0.003538518590750013
$ws.Range("J6").Value = 0.003538518590750013
$ws.Range("M6").Value = 1.819857
$ws.Range("N6").Value = 5.459571
$ws.Range("O6").Value = 0.01485317462584607
$ws.Range("P6").Value = 0.01485317462584607
$ws.Range("Q6").Value = 0.581602639059
$ws.Range("R6").Value = 5.234423751531001
$ws.Range("S6").Value = 0.00005255823454521268
$ws.Range("T6").Value = 0.00005255823454521268
# Row 7
$ws.Range("I7").Value = 0.003538518590750013
$ws.Range("J7").Value = 0.003538518590750013
$ws.Range("O7").Value = 0.726618572334523
$ws.Range("P7").Value = 0.7266185723345231
$ws.Range("S7").Value = 0.002571153326589943
$ws.Range("T7").Value = 0.002571153326589943
# Row 8
$ws.Range("I8").Value = 0.003538518590750013
$ws.Range("J8").Value = 0.003538518590750013
$ws.Range("M8").Value = 31.52924033333333
$ws.Range("N8").Value = 94.58772099999999
$ws.Range("O8").Value = 0.257333028084772
$ws.Range("P8").Value = 0.257333028084772
$ws.Range("Q8").Value = 10.076335330409
$ws.Range("R8").Value = 90.68701797368098
$ws.Range("S8").Value = 0.0009105777038919611
$ws.Range("T8").Value = 0.0009105777038919611
# Row 9
$ws.Range("I9").Value = 0.003538518590750013
$ws.Range("J9").Value = 0.003538518590750013
$ws.Range("K9").Value = 1.0
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.1464426666666667
$ws.Range("N9").Value = 0.439328
$ws.Range("O9").Value = 0.001195224954858853
$ws.Range("P9").Value = 0.001195224954858853
$ws.Range("Q9").Value = 0.046801172512
$ws.Range("R9").Value = 0.421210552608
$ws.Range("S9").Value = 0.000004229325722896395
$ws.Range("T9").Value = 0.000004229325722896396
# Row 10
$ws.Range("E10").Value = 1.0
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.115957
$ws.Range("H10").Value = 0.347871
$ws.Range("I10").Value = 0.001283894527085267
$ws.Range("J10").Value = 0.001283894527085267
$ws.Range("M10").Value = 1.819857
$ws.Range("N10").Value = 5.459571
$ws.Range("O10").Value = 0.01485317462584607
$ws.Range("P10").Value = 0.01485317462584607
$ws.Range("Q10").Value = 0.211025158149
$ws.Range("R10").Value = 1.899226423341
$ws.Range("S10").Value = 0.00001906990961196552
$ws.Range("T10").Value = 0.00001906990961196553
# Row 11
$ws.Range("E11").Value = 1.0
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.115957
$ws.Range("H11").Value = 0.347871
$ws.Range("I11").Value = 0.001283894527085267
$ws.Range("J11").Value = 0.001283894527085267
$ws.Range("O11").Value = 0.726618572334523
$ws.Range("P11").Value = 0.7266185723345231
$ws.Range("Q11").Value = 10.32336877492
$ws.Range("R11").Value = 92.91031897427999
$ws.Range("S11").Value = 0.0009329016082988042
$ws.Range("T11").Value = 0.0009329016082988045
# Row 12
$ws.Range("E12").Value = 1.0
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.115957
$ws.Range("H12").Value = 0.347871
$ws.Range("I12").Value = 0.001283894527085267
$ws.Range("J12").Value = 0.001283894527085267
$ws.Range("M12").Value = 31.52924033333333
$ws.Range("N12").Value = 94.58772099999999
$ws.Range("O12").Value = 0.257333028084772
$ws.Range("P12").Value = 0.257333028084772
$ws.Range("Q12").Value = 3.656036121332333
$ws.Range("R12").Value = 32.90432509199099
$ws.Range("S12").Value = 0.0003303884663963181
$ws.Range("T12").Value = 0.0003303884663963182
# Row 13
$ws.Range("E13").Value = 1.0
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.115957
$ws.Range("H13").Value = 0.347871
$ws.Range("I13").Value = 0.001283894527085267
$ws.Range("J13").Value = 0.001283894527085267
$ws.Range("K13").Value = 1.0
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.1464426666666667
$ws.Range("N13").Value = 0.439328
$ws.Range("O13").Value = 0.001195224954858853
$ws.Range("P13").Value = 0.001195224954858853
$ws.Range("Q13").Value = 0.01698105229866666
$ws.Range("R13").Value = 0.152829470688
$ws.Range("S13").Value = 0.000001534542778179016
$ws.Range("T13").Value = 0.000001534542778179017
# Row 14
$ws.Range("G14").Value = 89.26880233333334
$ws.Range("H14").Value = 267.806407
$ws.Range("I14").Value = 0.9883985163053822
$ws.Range("J14").Value = 0.9883985163053823
$ws.Range("M14").Value = 1.819857
$ws.Range("N14").Value = 5.459571
$ws.Range("O14").Value = 0.01485317462584607
$ws.Range("P14").Value = 0.01485317462584607
$ws.Range("Q14").Value = 162.456454807933
$ws.Range("R14").Value = 1462.108093271397
$ws.Range("S14").Value = 0.014680855762611
$ws.Range("T14").Value = 0.01468085576261101
# Row 15
$ws.Range("G15").Value = 89.26880233333334
$ws.Range("H15").Value = 267.806407
$ws.Range("I15").Value = 0.9883985163053822
$ws.Range("J15").Value = 0.9883985163053823
$ws.Range("O15").Value = 0.726618572334523
$ws.Range("P15").Value = 0.7266185723345231
$ws.Range("Q15").Value = 7947.383655858974
$ws.Range("R15").Value = 71526.45290273076
$ws.Range("S15").Value = 0.7181887188153776
$ws.Range("T15").Value = 0.7181887188153777
# Row 16
$ws.Range("G16").Value = 89.26880233333334
$ws.Range("H16").Value = 267.806407
$ws.Range("I16").Value = 0.9883985163053822
$ws.Range("J16").Value = 0.9883985163053823
$ws.Range("M16").Value = 31.52924033333333
$ws.Range("N16").Value = 94.58772099999999
$ws.Range("O16").Value = 0.257333028084772
$ws.Range("P16").Value = 0.257333028084772
$ws.Range("Q16").Value = 2814.577523036494
$ws.Range("R16").Value = 25331.19770732845
$ws.Range("S16").Value = 0.2543475831553599
$ws.Range("T16").Value = 0.25434758315536
# Row 17
$ws.Range("G17").Value = 89.26880233333334
$ws.Range("H17").Value = 267.806407
$ws.Range("I17").Value = 0.9883985163053822
$ws.Range("J17").Value = 0.9883985163053823
$ws.Range("K17").Value = 1.0
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.1464426666666667
$ws.Range("N17").Value = 0.439328
$ws.Range("O17").Value = 0.001195224954858853
$ws.Range("P17").Value = 0.001181358572033657
$ws.Range("Q17").Value = 13.07276146383289
$ws.Range("R17").Value = 117.654853174496
$ws.Range("S17").Value = 0.001181358572033657
$ws.Range("T17").Value = 0.001181358572033658
